# Added packing material script
$wb = $excel.ActiveWorkbook

# 1) Rename the "emp" sheet to "PackingMaterial"
$ws = $wb.Worksheets.Item("emp")
$ws.Name = "PackingMaterial"

# 2) Populate the sheet with the packing material data.
# The cells are written in the same order the original author typed them in
# (header row, then data row 2 left-to-right, then the remaining rows
# column-by-column) so that the shared-string table is built up in the
# same sequence as the canonical file.
$ws.Range("A1").Value = "MaterialName"
$ws.Range("B1").Value = "QTY"
$ws.Range("C1").Value = "Size"
$ws.Range("D1").Value = "Unit"
$ws.Range("E1").Value = "Cost"

$ws.Range("A2").Value = "a"
$ws.Range("B2").Value = 2
$ws.Range("C2").Value = "nos12"
$ws.Range("D2").Value = "pcs"
$ws.Range("E2").Value = 1

$ws.Range("A3").Value = "b"
$ws.Range("A4").Value = "v"
$ws.Range("A5").Value = "b"

$ws.Range("B3").Value = 3
$ws.Range("B4").Value = 2
$ws.Range("B5").Value = 3

$ws.Range("C3").Value = "3e"
$ws.Range("C4").Value = "ew"
$ws.Range("C5").Value = "we"

$ws.Range("D3").Value = "pcs"
$ws.Range("D4").Value = "pcs"
$ws.Range("D5").Value = "pcs"

$ws.Range("E3").Value = 1
$ws.Range("E4").Value = 1
$ws.Range("E5").Value = 1

# 3) Make this the active/selected sheet and cell
$ws.Range("E7").Select()
$ws.Activate()
